$wb = $excel.ActiveWorkbook

# Rename the first sheet from "PyToolConfig" to "Config"
$wsConfig = $wb.Worksheets.Item(1)
$wsConfig.Name = "Config"

# Make the Config sheet the active sheet/tab, and move its selection to E12
$wsConfig.Activate()
$wsConfig.Range("E12").Select()
